$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the RF (raising factor) column I for rows 22 through 41
# with the updated 2025 RF value.
$ws.Range("I22:I41").Value = 15.1448
